# Conserto do erro com o rotulo da coluna 2050 nas tabelas e retirada
# das linhas com total das tabelas.
#
# Sheets 1-3 and 5 have a year header row whose last column (E1) was
# mistakenly left as a leftover numeric value (676.8872865593887)
# instead of the "2050" text label used by the other year columns.
# Sheet 4 uses "2041-2050" (decade) as its label instead.
# Sheets 1-4 also have a trailing "Total" row (row 13) that must be
# removed, and sheet 6 has a trailing "Total" row (row 4) to remove.

$wb = $excel.ActiveWorkbook

function Set-YearLabel($ws, $cellRef, $label) {
    # A plain string assignment (Range.Value = "2050") auto-converts a
    # purely-numeric-looking string into a number cell, which is not
    # what we want (the other header cells are real text). Instead,
    # enter a text formula that evaluates to the label, then convert
    # it in place to a static value via copy / paste-values. This
    # keeps the cell's existing style untouched and avoids minting any
    # unused styles in styles.xml.
    $ws.Range($cellRef).Formula = '="' + $label + '"'
    $ws.Range($cellRef).Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)  # xlPasteValues
}

# --- Sheet 1: "Potencia Acumulada - SIN (MW)" ---
$ws1 = $wb.Worksheets.Item(1)
Set-YearLabel $ws1 "E1" "2050"
$ws1.Rows.Item(13).Delete()

# --- Sheet 2: "Geracao Periodo Medio (MWMed)" ---
$ws2 = $wb.Worksheets.Item(2)
Set-YearLabel $ws2 "E1" "2050"
$ws2.Rows.Item(13).Delete()

# --- Sheet 3: "Atendimento a Ponta(MW)" ---
$ws3 = $wb.Worksheets.Item(3)
Set-YearLabel $ws3 "E1" "2050"
$ws3.Rows.Item(13).Delete()

# --- Sheet 4: "Potencia Incremental - SIN(MW)" ---
# "2041-2050" is not purely numeric, so direct assignment keeps it text.
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("E1").Value = "2041-2050"
$ws4.Rows.Item(13).Delete()

# --- Sheet 5: "Emissoes Totais (MtCO2eq)" ---
$ws5 = $wb.Worksheets.Item(5)
Set-YearLabel $ws5 "E1" "2050"

# --- Sheet 6: "Custo Total (bilhoes de R$)" ---
$ws6 = $wb.Worksheets.Item(6)
$ws6.Rows.Item(4).Delete()

$excel.CutCopyMode = $false
